$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell $ws "B2" "Bitcoin"
Set-TextCell $ws "C2" "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextCell $ws "D2" "29.131.81"
Set-TextCell $ws "E2" "  -1.24%  "

Set-TextCell $ws "B3" "Ethereum"
Set-TextCell $ws "C3" "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextCell $ws "D3" "1.835.18"
Set-TextCell $ws "E3" "  -1.20%  "

Set-TextCell $ws "B4" "TetherUSD"
Set-TextCell $ws "C4" "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextCell $ws "D4" "0.9993"
Set-TextCell $ws "E4" "  -0.03%  "

Set-TextCell $ws "B5" "BNB"
Set-TextCell $ws "C5" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextCell $ws "D5" "240.14"
Set-TextCell $ws "E5" "  -2.14%  "

Set-TextCell $ws "B6" "XRP"
Set-TextCell $ws "C6" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell $ws "D6" "0.6647"
Set-TextCell $ws "E6" "  -4.56%  "

Set-TextCell $ws "B7" "USDC"
Set-TextCell $ws "C7" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextCell $ws "D7" "1.000"
Set-TextCell $ws "E7" "  +0.01%  "

Set-TextCell $ws "B8" "Cardano"
Set-TextCell $ws "C8" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell $ws "D8" "0.2954"
Set-TextCell $ws "E8" "  -4.04%  "

Set-TextCell $ws "B9" "Dogecoin"
Set-TextCell $ws "C9" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell $ws "D9" "0.07349"
Set-TextCell $ws "E9" "  -4.54%  "

Set-TextCell $ws "B10" "Solana"
Set-TextCell $ws "C10" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell $ws "D10" "22.73"
Set-TextCell $ws "E10" "  -3.97%  "

Set-TextCell $ws "B11" "TRON"
Set-TextCell $ws "C11" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws "D11" "0.07682"
Set-TextCell $ws "E11" "  -1.43%  "

Set-TextCell $ws "B12" "WrappedEther"
Set-TextCell $ws "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws "D12" "1.837.74"
Set-TextCell $ws "E12" "  -1.01%  "

Set-TextCell $ws "B13" "Polkadot"
Set-TextCell $ws "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws "D13" "5.019"
Set-TextCell $ws "E13" "  -2.75%  "

Set-TextCell $ws "B14" "Polygon"
Set-TextCell $ws "C14" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell $ws "D14" "0.6748"
Set-TextCell $ws "E14" "  -2.87%  "

Set-TextCell $ws "B15" "Litecoin"
Set-TextCell $ws "C15" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell $ws "D15" "86.12"
Set-TextCell $ws "E15" "  -5.61%  "

Set-TextCell $ws "B16" "Uniswap"
Set-TextCell $ws "C16" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell $ws "D16" "6.203"
Set-TextCell $ws "E16" "  -2.15%  "

Set-TextCell $ws "B17" "WrappedBTC"
Set-TextCell $ws "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell $ws "D17" "29.117.74"
Set-TextCell $ws "E17" "  -1.24%  "

Set-TextCell $ws "B18" "ShibaInu"
Set-TextCell $ws "C18" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell $ws "D18" "0.000008234"
Set-TextCell $ws "E18" "  -1.04%  "

Set-TextCell $ws "B19" "BitcoinCash"
Set-TextCell $ws "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell $ws "D19" "228.67"
Set-TextCell $ws "E19" "  -4.15%  "

Set-TextCell $ws "B20" "Avalanche"
Set-TextCell $ws "C20" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell $ws "D20" "12.50"
Set-TextCell $ws "E20" "  -1.94%  "

Set-TextCell $ws "B21" "Dai"
Set-TextCell $ws "C21" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws "D21" "0.9998"
Set-TextCell $ws "E21" "  -0.03%  "

Set-TextCell $ws "B22" "Chainlink"
Set-TextCell $ws "C22" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws "D22" "7.301"
Set-TextCell $ws "E22" "  -4.27%  "

Set-TextCell $ws "B23" "BinanceUSD"
Set-TextCell $ws "C23" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell $ws "D23" "1.000"
Set-TextCell $ws "E23" "  +0.01%  "

Set-TextCell $ws "B24" "Monero"
Set-TextCell $ws "C24" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D24" "161.09"
Set-TextCell $ws "E24" "  +0.62%  "

Set-TextCell $ws "B25" "Stellar"
Set-TextCell $ws "C25" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D25" "0.1417"
Set-TextCell $ws "E25" "  -5.28%  "

Set-TextCell $ws "B26" "Cosmos"
Set-TextCell $ws "C26" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws "D26" "8.673"
Set-TextCell $ws "E26" "  -2.57%  "

Set-TextCell $ws "B27" "EthereumClassic"
Set-TextCell $ws "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws "D27" "18.03"
Set-TextCell $ws "E27" "  -1.36%  "

Set-TextCell $ws "B28" "PancakeSwap"
Set-TextCell $ws "C28" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws "D28" "1.502"
Set-TextCell $ws "E28" "  -2.20%  "

Set-TextCell $ws "B29" "Filecoin"
Set-TextCell $ws "C29" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws "D29" "4.229"
Set-TextCell $ws "E29" "  -0.51%  "

Set-TextCell $ws "B30" "InternetComputer(DFINITY)"
Set-TextCell $ws "C30" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws "D30" "4.096"
Set-TextCell $ws "E30" "  -1.28%  "

Set-TextCell $ws "B31" "Toncoin"
Set-TextCell $ws "C31" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws "D31" "1.199"
Set-TextCell $ws "E31" "  -0.50%  "

Set-TextCell $ws "B32" "Hedera"
Set-TextCell $ws "C32" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws "D32" "0.05304"
Set-TextCell $ws "E32" "  +3.80%  "

Set-TextCell $ws "B33" "LidoDAOToken"
Set-TextCell $ws "C33" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws "D33" "1.858"
Set-TextCell $ws "E33" "  -1.27%  "

Set-TextCell $ws "B34" "ImmutableX"
Set-TextCell $ws "C34" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws "D34" "0.7468"
Set-TextCell $ws "E34" "  -3.78%  "

Set-TextCell $ws "B35" "ARBITRUM"
Set-TextCell $ws "C35" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws "D35" "1.129"
Set-TextCell $ws "E35" "  -1.71%  "

Set-TextCell $ws "B36" "HuobiToken"
Set-TextCell $ws "C36" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell $ws "D36" "2.681"
Set-TextCell $ws "E36" "  -0.18%  "

Set-TextCell $ws "B37" "Maker"
Set-TextCell $ws "C37" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws "D37" "1.316.98"
Set-TextCell $ws "E37" "  -0.03%  "

Set-TextCell $ws "B38" "VeChain"
Set-TextCell $ws "C38" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D38" "0.01806"
Set-TextCell $ws "E38" "  -3.80%  "

Set-TextCell $ws "B39" "MXToken"
Set-TextCell $ws "C39" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell $ws "D39" "2.715"
Set-TextCell $ws "E39" "  -0.41%  "

Set-TextCell $ws "B40" "TrustWalletToken"
Set-TextCell $ws "C40" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws "D40" "0.9224"
Set-TextCell $ws "E40" "  -3.19%  "

Set-TextCell $ws "B41" "FraxShare"
Set-TextCell $ws "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws "D41" "5.980"
Set-TextCell $ws "E41" "  +3.62%  "

Set-TextCell $ws "B42" "PaxDollar"
Set-TextCell $ws "C42" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell $ws "D42" "0.9983"
Set-TextCell $ws "E42" "  -0.27%  "

Set-TextCell $ws "B43" "Quant"
Set-TextCell $ws "C43" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell $ws "D43" "103.41"
Set-TextCell $ws "E43" "  -2.48%  "

Set-TextCell $ws "B44" "RocketPoolETH"
Set-TextCell $ws "C44" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell $ws "D44" "1.984.62"
Set-TextCell $ws "E44" "  -0.82%  "

Set-TextCell $ws "B45" "Mantle"
Set-TextCell $ws "C45" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws "D45" "0.5169"
Set-TextCell $ws "E45" "  -1.26%  "

Set-TextCell $ws "B46" "BabyDogeCoin"
Set-TextCell $ws "C46" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell $ws "D46" "0.00000000121"
Set-TextCell $ws "E46" "  -3.22%  "

Set-TextCell $ws "B47" "Aave"
Set-TextCell $ws "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws "D47" "63.82"
Set-TextCell $ws "E47" "  +1.00%  "

Set-TextCell $ws "B48" "RenderToken"
Set-TextCell $ws "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws "D48" "1.760"
Set-TextCell $ws "E48" "  -1.57%  "

Set-TextCell $ws "B49" "EnergySwap"
Set-TextCell $ws "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D49" "9.256"
Set-TextCell $ws "E49" "  -5.64%  "

Set-TextCell $ws "B50" "XinFinNetwork"
Set-TextCell $ws "C50" "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
Set-TextCell $ws "D50" "0.07420"
Set-TextCell $ws "E50" "  +8.87%  "

Set-TextCell $ws "B51" "Cronos"
Set-TextCell $ws "C51" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws "D51" "0.05933"
Set-TextCell $ws "E51" "  +0.06%  "

